$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update NIK column (C2:C7) from "EN-4-047" to "EN-4-046"
$ws.Range("C2:C7").Value = "EN-4-046"

# Update Nama Karyawan column (D2:D7) from "Ari Pratama " to "Agus Priyanto"
$ws.Range("D2:D7").Value = "Agus Priyanto"
